$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.175.82"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.052.86"
$ws.Range("E3").Value = "  +3.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.07"
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("E7").Value = "  +0.01%  "
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.36"
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = "  +5.07%  "
$ws.Range("E9").Value = "  +2.56%  "
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.93"
$ws.Range("D10").Style = $style
$ws.Range("E10").Value = "  -1.99%  "
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("D13").Value = "2.357.33"
$ws.Range("E13").Value = "  +4.05%  "
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.47"
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = "  +4.46%  "
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.75"
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = "  -0.02%  "
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.771"
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = "  +3.31%  "
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.23"
$ws.Range("D17").Style = $style
$ws.Range("E17").Value = "  +4.90%  "
$ws.Range("D18").Value = "2.069.16"
$ws.Range("E18").Value = "  +4.93%  "
$ws.Range("D19").Value = "37.333.16"
$ws.Range("E19").Value = "  +1.39%  "
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.91"
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = "  +20.26%  "
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.30"
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "0.0₃0805"
$ws.Range("E22").Value = "  +0.33%  "
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "223.01"
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = "  -1.54%  "
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = "  +3.57%  "
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.74"
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = "  +0.92%  "
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.84"
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = "  +3.05%  "
$ws.Range("E29").Value = "  +5.97%  "
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.15"
$ws.Range("D30").Style = $style
$ws.Range("E30").Value = "  +1.16%  "
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.36"
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = "  +6.60%  "
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.118"
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = "  +1.34%  "
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.45"
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("E35").Value = "  +9.62%  "
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.35"
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = "  +3.12%  "
$ws.Range("E37").Value = "  -0.06%  "
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.92"
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = "  +15.94%  "
$ws.Range("E39").Value = "  +1.41%  "
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = "  -1.14%  "
$ws.Range("E41").Value = "  -1.90%  "
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0949"
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = "  +7.07%  "
$ws.Range("D43").Value = "1.477.27"
$ws.Range("E43").Value = "  +4.54%  "
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.28"
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = "  +11.33%  "
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.18"
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = "  +8.60%  "
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("E47").Value = "  +1.36%  "
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.95"
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = "  +5.63%  "
$ws.Range("E49").Value = "  +3.01%  "
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.15"
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = "  +8.89%  "
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.93"
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = "  +2.54%  "
